# "Adding fillers to the test"
#
# The header row labelled the four math-operation columns simply
# "addition" / "subtraction" / "multiplication" / "division". To make room
# for a new "filler" item type (a plain sentence-completion item, as
# opposed to the math items), the author renamed those four headers to be
# explicitly prefixed with "math_" so they read as
# "math_addition" / "math_subtraction" / "math_multiplication" / "math_division".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "math_addition"
$ws.Range("D1").Value = "math_subtraction"
$ws.Range("E1").Value = "math_multiplication"
$ws.Range("F1").Value = "math_division"

# Leave the cursor where the author last left it before saving.
[void]$ws.Range("F2").Select()
